$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 456.19232
$ws.Range("I33").Value = 87.47369
$ws.Range("J33").Value = 1457
$ws.Range("K33").Value = 87.47369
$ws.Range("L33").Value = 1457
$ws.Range("M33").Value = 141.52631
$ws.Range("N33").Value = -1915
$ws.Range("H112").Value = 52633656
$ws.Range("J112").Value = 2365.3125
$ws.Range("L112").Value = 7095.9375
$ws.Range("N112").Value = -9311.9375
$ws.Range("H113").Value = 2756.1667
$ws.Range("I113").Value = 2533.2222
$ws.Range("J113").Value = 3157.4666
$ws.Range("K113").Value = 2533.2222
$ws.Range("L113").Value = 3157.4666
$ws.Range("M113").Value = 720.7777999999998
$ws.Range("N113").Value = -9665.4666
$ws.Range("H116").Value = 1836.909
$ws.Range("J116").Value = 2276.5
$ws.Range("L116").Value = 2276.5
$ws.Range("N116").Value = -9160.5
$ws.Range("H129").Value = 860.375
$ws.Range("I129").Value = 295.2
$ws.Range("J129").Value = 1117.2727
$ws.Range("K129").Value = 885.5999999999999
$ws.Range("L129").Value = 3351.8181
$ws.Range("M129").Value = 4114.4
$ws.Range("N129").Value = -13351.8181
$ws.Range("H132").Value = 832919.1
$ws.Range("I132").Value = 2696.1353
$ws.Range("J132").Value = 2229203.2
$ws.Range("K132").Value = 8088.4059
$ws.Range("L132").Value = 6687609.600000001
$ws.Range("M132").Value = -5558.4059
$ws.Range("N132").Value = -6692669.600000001
$ws.Range("H137").Value = 1755848.6
$ws.Range("I137").Value = 2440026.2
$ws.Range("K137").Value = 7320078.600000001
$ws.Range("M137").Value = -7317528.600000001
$ws.Range("H138").Value = 2138719.5
$ws.Range("I138").Value = 1373.3478
$ws.Range("J138").Value = 5211154.5
$ws.Range("K138").Value = 4120.0434
$ws.Range("L138").Value = 15633463.5
$ws.Range("M138").Value = 1019.9566
$ws.Range("N138").Value = -15643743.5
$ws.Range("H141").Value = 2729.6223
$ws.Range("I141").Value = 1561.1316
$ws.Range("K141").Value = 4683.3948
$ws.Range("M141").Value = 496.6052

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3627.25
$ws.Range("I32").Value = 2864.9639
$ws.Range("J32").Value = 7349
$ws.Range("K32").Value = 2864.9639
$ws.Range("L32").Value = 7349
$ws.Range("M32").Value = -2577.9639
$ws.Range("N32").Value = -7923
$ws.Range("H61").Value = 18219572
$ws.Range("I61").Value = 20855272
$ws.Range("J61").Value = 146202
$ws.Range("K61").Value = 20855272
$ws.Range("L61").Value = 146202
$ws.Range("M61").Value = -20855060
$ws.Range("N61").Value = -146626
$ws.Range("H132").Value = 62687.117
$ws.Range("I132").Value = 40861.12
$ws.Range("J132").Value = 123314.89
$ws.Range("K132").Value = 122583.36
$ws.Range("L132").Value = 369944.67
$ws.Range("M132").Value = -120053.36
$ws.Range("N132").Value = -375004.67
$ws.Range("H136").Value = 18219572
$ws.Range("I136").Value = 20855272
$ws.Range("J136").Value = 146202
$ws.Range("K136").Value = 62565816
$ws.Range("L136").Value = 438606
$ws.Range("M136").Value = -62563266
$ws.Range("N136").Value = -443706

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").Value = ""
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").Value = ""
$ws.Range("H134").Value = 1795.0817
$ws.Range("I134").Value = 1159.0857
$ws.Range("K134").Value = 3477.2571
$ws.Range("M134").Value = -942.2571000000003

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 13335092
$ws.Range("I16").Value = 1414.1666
$ws.Range("K16").Value = 1414.1666
$ws.Range("M16").Value = -1127.1666
$ws.Range("H45").Value = 17500
$ws.Range("I45").Value = 15000
$ws.Range("J45").Value = 20000
$ws.Range("K45").Value = 15000
$ws.Range("L45").Value = 20000
$ws.Range("M45").Value = -14407
$ws.Range("N45").Value = -21186
$ws.Range("H58").Value = 20835080
$ws.Range("I58").Value = 27028626
$ws.Range("K58").Value = 27028626
$ws.Range("M58").Value = -27028423
$ws.Range("H99").Value = 5683.3335
$ws.Range("I99").Value = 2975
$ws.Range("J99").Value = 11100
$ws.Range("K99").Value = 2975
$ws.Range("L99").Value = 11100
$ws.Range("M99").Value = -1477
$ws.Range("N99").Value = -14096
$ws.Range("H113").Value = 13335092
$ws.Range("I113").Value = 1414.1666
$ws.Range("K113").Value = 1414.1666
$ws.Range("M113").Value = 755.8334
$ws.Range("H126").Value = 5683.3335
$ws.Range("I126").Value = 2975
$ws.Range("J126").Value = 11100
$ws.Range("K126").Value = 8925
$ws.Range("L126").Value = 33300
$ws.Range("M126").Value = -6455
$ws.Range("N126").Value = -38240
$ws.Range("H132").Value = 16484.537
$ws.Range("I132").Value = 1372.8368
$ws.Range("K132").Value = 4118.5104
$ws.Range("M132").Value = -1588.5104
$ws.Range("H134").Value = 18459.984
$ws.Range("I134").Value = 1297.3914
$ws.Range("K134").Value = 3892.1742
$ws.Range("M134").Value = -1357.1742
$ws.Range("H135").Value = 40180.6
$ws.Range("J135").Value = 41714.285
$ws.Range("L135").Value = 41714.285
$ws.Range("N135").Value = -51854.285
$ws.Range("H136").Value = 20835080
$ws.Range("I136").Value = 27028626
$ws.Range("K136").Value = 81085878
$ws.Range("M136").Value = -81083328

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 2670.875
$ws.Range("J106").Value = 2670.875
$ws.Range("L106").Value = 8012.625
$ws.Range("N106").Value = -9904.625

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1919.3182
$ws.Range("I126").Value = 1167.6666
$ws.Range("J126").Value = 2821.3
$ws.Range("K126").Value = 3502.9998
$ws.Range("L126").Value = 8463.900000000001
$ws.Range("M126").Value = -1032.9998
$ws.Range("N126").Value = -13403.9
$ws.Range("H132").Value = 38499.906
$ws.Range("I132").Value = 30707.5
$ws.Range("J132").Value = 51747
$ws.Range("K132").Value = 92122.5
$ws.Range("L132").Value = 155241
$ws.Range("M132").Value = -89592.5
$ws.Range("N132").Value = -160301

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 20858.79
$ws.Range("I132").Value = 1305.9032
$ws.Range("J132").Value = 49722.57
$ws.Range("K132").Value = 3917.7096
$ws.Range("L132").Value = 149167.71
$ws.Range("M132").Value = -1387.7096
$ws.Range("N132").Value = -154227.71

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").Value = ""
$ws.Range("H132").Value = 38065.85
$ws.Range("I132").Value = 25695.45
$ws.Range("J132").Value = 73409.86
$ws.Range("K132").Value = 77086.35000000001
$ws.Range("L132").Value = 220229.58
$ws.Range("M132").Value = -74556.35000000001
$ws.Range("N132").Value = -225289.58
